$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: average of the |S*|/n column (J)
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"

# Rows 14-17: summary statistics (labels in column A, values/formulas in column B)
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"

$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"

$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"

$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"

# Format B14 with bold, 12pt, vertically centered font, then copy that
# formatting to B15:B17 so all four cells share a single new cell style
# (mirrors the "apply once, paint elsewhere" workflow of the original author).
$r14 = $ws.Range("B14")
$r14.Font.Bold = $true
$r14.Font.Size = 12
$r14.VerticalAlignment = -4108   # xlCenter

$r14.Copy()
$ws.Range("B15:B17").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# Row heights for the new summary rows
$ws.Rows.Item(14).RowHeight = 15.6
$ws.Rows.Item(15).RowHeight = 15.6
$ws.Rows.Item(16).RowHeight = 15.6
$ws.Rows.Item(17).RowHeight = 15.6

# Page setup (paper size / orientation)
$ws.PageSetup.PaperSize = 9      # xlPaperA4
$ws.PageSetup.Orientation = 1    # xlPortrait

# Selection matches the saved view in the target workbook
$ws.Range("A14:B17").Select()
